{"js": "// 1) The \"_GoBack\" bookmark currently sits alone in the empty paragraph\n//    right after the title. Remove it from there \u2014 that paragraph becomes\n//    a plain empty paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the list paragraph that holds\n//    \"Verificar las actividades hecha por los integrantes \".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Verificar las actividades hecha por los integrantes\") !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error('Could not find anchor paragraph \"Verificar las actividades hecha por los integrantes\".');\n}\n\n// 3) Insert a brand-new list paragraph right after it, seeded with the\n//    first run's text (it inherits the list/numbering + run formatting\n//    from the anchor paragraph automatically).\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"Ver la tabla de actividades y su poca duraci\u00f3n\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 4) Append the trailing run \" en el GANTT \" at the end of the new\n//    paragraph (still a plain text append for now).\nconst trailingRunRange = newParagraph.insertText(\" en el GANTT \", Word.InsertLocation.end);\nawait context.sync();\n\n// 5) Re-insert the \"_GoBack\" bookmark exactly at the boundary between the\n//    two runs, i.e. right before the trailing \" en el GANTT \" text. This\n//    splits the paragraph text into two separate runs with the bookmark\n//    sandwiched between them, matching the target structure.\nconst trailingRunStart = trailingRunRange.getRange(Word.RangeLocation.start);\ntrailingRunStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The \"_GoBack\" bookmark currently sits alone in the empty paragraph\n#    right after the title. Remove it from there \u2014 that paragraph becomes\n#    a plain empty paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the list paragraph that holds\n#    \"Verificar las actividades hecha por los integrantes \" and insert a\n#    brand-new list paragraph right after it (it inherits the list /\n#    numbering + run formatting from the anchor paragraph automatically).\n$findRange = $d.Content\n$findRange.Find.Execute(\"Verificar las actividades hecha por los integrantes\")\n$anchorPara = $findRange.Paragraphs(1)\n$anchorPara.Range.InsertParagraphAfter()\n\n$newPara = $anchorPara.Next()\n$newRange = $newPara.Range\n\n$firstSentence = \"Ver la tabla de actividades y su poca duraci\u00f3n\"\n$secondSentence = \" en el GANTT \"\n\n# Insert the whole new sentence in one shot first (inserting the two\n# pieces as one string avoids ending up with a bookmark sitting on the\n# very last character of the paragraph, right before the paragraph\n# mark, which this engine otherwise snaps to span the whole paragraph).\n$newRange.InsertBefore($firstSentence + $secondSentence)\n\n# 3) Re-insert the \"_GoBack\" bookmark exactly at the boundary between the\n#    two sentences, splitting the paragraph text into two separate runs\n#    with the bookmark sandwiched between them.\n$splitPos = $newRange.Start + $firstSentence.Length\n$bmRange = $d.Range($splitPos, $splitPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
